$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'0304"
$ws.Range("C2").Value = "Cimentão"
$ws.Range("D2").Value = "DVR 2 está sem comunicação de câmeras, era via DDNS e pedi ele pra manter."
$ws.Range("E2").Value = "Técnico disse que vai no dentista, pode ser que ele faça o atendimento ainda hoje. Disse que amanhã vai compensar nas o.s."
$ws.Range("G2").Value = "Pendente"

$ws.Range("B3").Value = "'0563"
$ws.Range("C3").Value = "Unicred"
$ws.Range("D3").Value = "Cliente pedindo manutenção em câmeras não monitoradas."
$ws.Range("E3").Value = ""
$ws.Range("G3").Value = "Em andamento"

$ws.Range("A4").Value = "Ryan"
$ws.Range("B4").Value = "'0870"
$ws.Range("C4").Value = "Colégio Santo Agostinho"
$ws.Range("D4").Value = "Local tem duas zonas abertas, cliente tem uma AMT 8000. Passei 4 bateria maior e 4 bateria redondinha pro técnico conferir isso pra nós."
$ws.Range("E4").Value = ""
$ws.Range("G4").Value = "Em andamento"

$ws.Range("A5:G10").ClearContents()

$ws.Rows("3:3").AutoFit()
$ws.Rows("5:10").AutoFit()

$ws.Range("H4").Select()
